$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.081.12'
$ws.Range('D3').Value = '1.651.29'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '217.45'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').Value = '0.5259'
$ws.Range('E6').Value = '  +2.37%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '0.2595'
$ws.Range('E8').Value = '  -1.87%  '
$ws.Range('D9').Value = '0.06331'
$ws.Range('E9').Value = '  +0.94%  '
$ws.Range('E10').Value = '  -2.25%  '
$ws.Range('D11').Value = '0.07802'
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('D12').Value = '4.505'
$ws.Range('E12').Value = '  +0.89%  '
$ws.Range('D13').Value = '1.675.44'
$ws.Range('E13').Value = '  +0.98%  '
$ws.Range('D14').Value = '0.5485'
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('D15').Value = '0.0₅8220'
$ws.Range('E15').Value = '  +1.46%  '
$ws.Range('D16').Value = '65.34'
$ws.Range('E16').Value = '  +0.72%  '
$ws.Range('D17').Value = '26.088.23'
$ws.Range('E17').Value = '  -0.42%  '
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D19').Value = '4.571'
$ws.Range('E19').Value = '  -1.02%  '
$ws.Range('D20').Value = '191.11'
$ws.Range('E20').Value = '  -0.50%  '
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('D22').Value = '6.030'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('E23').Value = '  -0.11%  '
$ws.Range('D24').Value = '142.30'
$ws.Range('E24').Value = '  +1.63%  '
$ws.Range('D25').Value = '0.1236'
$ws.Range('E25').Value = '  +1.27%  '
$ws.Range('D26').Value = '7.225'
$ws.Range('E26').Value = '  -0.47%  '
$ws.Range('E27').Value = '  -0.65%  '
$ws.Range('D28').Value = '1.429'
$ws.Range('E28').Value = '  -0.45%  '
$ws.Range('D29').Value = '0.05811'
$ws.Range('E29').Value = '  -2.51%  '
$ws.Range('E30').Value = '  +0.15%  '
$ws.Range('E31').Value = '  -0.74%  '
$ws.Range('D32').Value = '3.257'
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('D33').Value = '1.584'
$ws.Range('E33').Value = '  -0.72%  '
$ws.Range('D34').Value = '2.412'
$ws.Range('D35').Value = '0.9442'
$ws.Range('E35').Value = '  -2.29%  '
$ws.Range('E36').Value = '  +0.22%  '
$ws.Range('D37').Value = '0.5724'
$ws.Range('E37').Value = '  +0.45%  '
$ws.Range('E38').Value = '  +0.88%  '
$ws.Range('D39').Value = '0.8444'
$ws.Range('E39').Value = '  -1.67%  '
$ws.Range('D40').Value = '5.746'
$ws.Range('E40').Value = '  -4.15%  '
$ws.Range('D41').Value = '1.002'
$ws.Range('E42').Value = '  +3.08%  '
$ws.Range('D43').Value = '1.026.14'
$ws.Range('E43').Value = '  +1.54%  '
$ws.Range('D44').Value = '1.796.44'
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('D45').Value = '56.97'
$ws.Range('E45').Value = '  +0.47%  '
$ws.Range('D46').Value = '0.9992'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('D48').Value = '0.05146'
$ws.Range('E48').Value = '  -0.39%  '
$ws.Range('E49').Value = '  +0.64%  '
$ws.Range('D50').Value = '7.816'
$ws.Range('E50').Value = '  -2.71%  '
$ws.Range('D51').Value = '0.09642'
